# Normalize the "Recorded By" (column G) entries so that any variant of
# "System" that appears as the second (or later) name in a comma-separated
# list of recorders is moved to the front - i.e. the first two
# comma-separated entries are swapped whenever the value does not already
# start with the exact token "System".
#
# Examples:
#   "dnasr281@gmail.com, System"                 -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"                     -> "System, admin@admin.com"
#   "system, System, backup@backdoor.com"         -> "System, system, backup@backdoor.com"
#   "admin@admin.com, dnasr281@gmail.com"         -> "dnasr281@gmail.com, admin@admin.com"
#   "System, backup@backdoor.com"                 -> unchanged (already starts with "System")
#   "dnasr281@gmail.com"                          -> unchanged (single value, no comma)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G" + $r)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split ",\s*"
    if ($parts.Count -lt 2) { continue }
    if ($parts[0].Equals("System")) { continue }

    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp

    $cell.Value2 = [string]::Join(", ", $parts)
}
